# cadastro_empresas fix: append the two missing employee records
# (Paulo / Git e Github) to the "empregados" sheet.
#
# Sheet layout (row 1 = header): Nome | Email | Empresa
#   Row 2 (existing): Guilherme | celente.guilherme@outlook.com | rtRfd34
#   Row 3 (new):      Paulo     | progeri@yahoo.com             | rtRfd34
#   Row 4 (new):      Git e Github | guilherme.celente@escola.pr.gov.br | rtRfd34

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New employee: Paulo
$ws.Range("A3").Value = "Paulo"
$ws.Range("B3").Value = "progeri@yahoo.com"
$ws.Range("C3").Value = "rtRfd34"

# New employee: Git e Github
$ws.Range("A4").Value = "Git e Github"
$ws.Range("B4").Value = "guilherme.celente@escola.pr.gov.br"
$ws.Range("C4").Value = "rtRfd34"
